# Deploy the implementation guide: refresh the generated Status/Date
# metadata rows on the "Metadata" sheet of the FHIR ValueSet spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$newStatus = "draft"
$newDate   = "2023-08-01T16:12:28+00:00"

# Find the "Status" / "Date" property rows in column A (labels live in A,
# values in B) instead of hard-coding row numbers, and update the value
# cell next to each label.
$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $label = $ws.Cells.Item($r, 1).Value()

    if ($label -eq "Status") {
        $ws.Cells.Item($r, 2).Value = $newStatus
    }
    elseif ($label -eq "Date") {
        $ws.Cells.Item($r, 2).Value = $newDate
    }
}
